# Week 13 logging update
# Updates player stat totals on the "Rushing" and "Receiving" sheets to
# reflect the latest week's accumulated numbers.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$wsRushing = $wb.Worksheets.Item("Rushing")

# Row 2 - J.Allen
$wsRushing.Range("C2").Value = 19
$wsRushing.Range("D2").Value = 23
$wsRushing.Range("E2").Value = 27
$wsRushing.Range("F2").Value = 19

# Row 3 - D.Singletary
$wsRushing.Range("C3").Value = 56
$wsRushing.Range("D3").Value = 44
$wsRushing.Range("F3").Value = 14

# Row 4 - Z.Moss
$wsRushing.Range("C4").Value = 42
$wsRushing.Range("D4").Value = 28
$wsRushing.Range("E4").Value = 8
$wsRushing.Range("F4").Value = 25

# Row 5 - M.Breida
$wsRushing.Range("C5").Value = 11
$wsRushing.Range("D5").Value = 9

# --- Receiving sheet ---
$wsReceiving = $wb.Worksheets.Item("Receiving")

# Row 3 - Z.Moss
$wsReceiving.Range("C3").Value = 28
$wsReceiving.Range("D3").Value = 19
$wsReceiving.Range("G3").Value = 6
$wsReceiving.Range("H3").Value = 4

# Row 4 - M.Breida
$wsReceiving.Range("C4").Value = 7
$wsReceiving.Range("D4").Value = 6

# Row 6 - S.Diggs
$wsReceiving.Range("C6").Value = 84
$wsReceiving.Range("D6").Value = 61
$wsReceiving.Range("E6").Value = 25
$wsReceiving.Range("F6").Value = 10
$wsReceiving.Range("G6").Value = 20

# Row 7 - E.Sanders
$wsReceiving.Range("C7").Value = 45
$wsReceiving.Range("D7").Value = 31
$wsReceiving.Range("E7").Value = 26

# Row 8 - C.Beasley
$wsReceiving.Range("C8").Value = 74
$wsReceiving.Range("D8").Value = 58
$wsReceiving.Range("G8").Value = 12

# Row 9 - G.Davis
$wsReceiving.Range("C9").Value = 16
$wsReceiving.Range("D9").Value = 9
$wsReceiving.Range("E9").Value = 13
$wsReceiving.Range("G9").Value = 8
$wsReceiving.Range("H9").Value = 4

# Row 12 - D.Knox
$wsReceiving.Range("C12").Value = 33
$wsReceiving.Range("D12").Value = 27
$wsReceiving.Range("E12").Value = 14
$wsReceiving.Range("G12").Value = 11
